# Update the verified dashboard metrics for shop-insight-v2
# Append 7 new email rows (80-86) to Sheet1, column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newEmails = @(
    "fcn07453@zslsz.com",
    "zlv91813@zslsz.com",
    "ell61810@zbock.com",
    "xtz13936@nezid.com",
    "ufe42269@zbock.com",
    "bux65944@zbock.com",
    "rqb91177@zbock.com"
)

$startRow = 80
for ($i = 0; $i -lt $newEmails.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newEmails[$i]
}
